$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the 2nd and 3rd row of every 4-row "year" block (rows 2-5, 6-9, ... 62-65).
# i.e. swap rows (3,4), (7,8), (11,12), ... (63,64) across columns A:E.
for ($start = 2; $start -le 62; $start += 4) {
    $row1 = $start + 1
    $row2 = $start + 2

    $r1 = $ws.Range("A$row1" + ":E$row1")
    $r2 = $ws.Range("A$row2" + ":E$row2")

    $v1 = $r1.Value2
    $v2 = $r2.Value2

    $r1.Value2 = $v2
    $r2.Value2 = $v1
}

# Remove columns F and G entirely (data no longer tracked).
$ws.Range("F1:G1").EntireColumn.Delete()
